$d = $word.ActiveDocument

$replacements = @(
    @("345÷8=", "640÷8="),
    @("224÷8=", "656÷9="),
    @("873÷3=", "771÷2="),
    @("611÷9=", "996÷2="),
    @("639÷9=", "352÷2="),
    @("563÷6=", "128÷8="),
    @("200÷5=", "113÷3="),
    @("365÷2=", "506÷3="),
    @("547÷8=", "537÷4="),
    @("865÷5=", "290÷4="),
    @("377÷4=", "455÷6="),
    @("716÷4=", "816÷5="),
    @("740÷6=", "363÷2="),
    @("711÷6=", "304÷5="),
    @("363÷6=", "967÷3="),
    @("657÷2=", "864÷4="),
    @("995÷5=", "143÷9="),
    @("560÷8=", "613÷6="),
    @("768÷7=", "537÷3="),
    @("239÷5=", "695÷9="),
    @("972÷4=", "342÷9="),
    @("932÷2=", "119÷9="),
    @("479÷7=", "434÷8="),
    @("965÷8=", "645÷5="),
    @("568÷7=", "695÷3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
